$d = $word.ActiveDocument

# --- 1. Change the letter date from "12" to "13" (August 12, 2016 -> August 13, 2016) ---
$dateRange = $d.Content
$dateRange.Find.ClearFormatting()
$found = $dateRange.Find.Execute("12", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateRange.Text = "13"

# --- 2. Move the automatic "_GoBack" bookmark to mark this as the most recent edit,
#        exactly like Word does when you edit and then save the document. ---
$newSpot = $d.Range($dateRange.End, $dateRange.End)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $newSpot)

# --- 3. The old "_GoBack" bookmark used to sit between "net" and "work exchange formats...";
#        now that it has moved away, merge those two runs back into a single run of text
#        (matching how Word recombines the split once the bookmark no longer separates them). ---
$netRange = $d.Content
$netRange.Find.ClearFormatting()
$netFound = $netRange.Find.Execute("commonly used net", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$netRange.Collapse(0)
$junction = $netRange.Start
$netRange.InsertAfter("X")
$d.Range($junction, $junction + 1).Text = ""
